# modified test cases on overdue fix
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: recompute the overdue-fee row (row 5) and drop the now
# -unused trailing zero row (row 6); also touches G2 / selection.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A5").Value = 17.76
$wsSummary.Range("E5").Value = 17.76
$wsSummary.Range("F5").Value = 17.76
$wsSummary.Cells.Item(2, 7).Style = "Normal"
$wsSummary.Rows.Item(6).Delete()
$wsSummary.Range("D5").Select()

# ---------------------------------------------------------------------
# Repayment schedule sheet: the Over Due column (O) collapses into the
# sheet - its zero placeholders are dropped - and the 3rd instalment
# (row 5) is recalculated now that the fee no longer applies.
# ---------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Range("P2").Clear()
$wsSchedule.Range("O3").Clear()
$wsSchedule.Range("O4").Clear()
$wsSchedule.Range("J5").Value = 0
$wsSchedule.Range("K5").Value = 887.72
$wsSchedule.Range("O5").Clear()
$wsSchedule.Range("P5").Value = 887.72
$wsSchedule.Range("O6").Clear()
$wsSchedule.Range("O7").Clear()
$wsSchedule.Range("O8").Clear()
$wsSchedule.Range("G6").Select()

# ---------------------------------------------------------------------
# Transactions sheet becomes the active tab / selection, NewLoanInput
# loses it.
# ---------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("D2").Select()
